$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mets")
$ws.Range("F1").Value = "measured?"
$ws.Range("F1").Font.Bold = $true
$ws.Range("F2:F21").Value = 0
